$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder the "Recorded By" e-mail list in G2 ---
$ws.Range("G2").Value = "gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# --- Widen column I (9) from 10 to 14 characters ---
# ColumnWidth is stored with a +5/6 padding and quantized to 1/6ths by the
# engine, so 13.14 round-trips to an on-disk width of exactly 14.
$ws.Columns.Item(9).ColumnWidth = 13.14

# --- Missing/Pending session counters ---
$ws.Range("L7").Value = 1
$ws.Range("L8").Value = 26
$ws.Range("P15").Value = 1
$ws.Range("Q15").Value = 26

# --- Re-style row 28 (new PHYSIOLOGY session) to the "Not Recorded" look:
#     black text on a pink fill, centered - and update its status text.
$row28 = $ws.Range("A28:I28")
$row28.Font.Color = 0
$row28.Interior.Color = 12695295
$row28.HorizontalAlignment = -4108
$row28.VerticalAlignment = -4108

$ws.Range("I28").Value = "Not Recorded"
